$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptocurrency "Price" (column D) and "Volume(1h)" (column E)
# values, refreshed by the scheduled scraper run.
$cellUpdates = @{
    "D2" = "69.898.50"
    "E2" = "  +3.89%  "
    "D3" = "3.585.78"
    "E3" = "  +3.81%  "
    "E4" = "  +0.09%  "
    "D5" = "585.95"
    "E5" = "  +3.08%  "
    "D6" = "189.06"
    "E6" = "  +2.63%  "
    "D7" = "0.640"
    "E7" = "  +1.08%  "
    "D8" = "3.578.24"
    "E8" = "  +3.80%  "
    "E9" = "  +0.07%  "
    "E10" = "  -0.65%  "
    "D11" = "0.656"
    "E11" = "  +1.67%  "
    "D12" = "57.54"
    "E12" = "  +3.52%  "
    "D13" = "0.0000287"
    "E13" = "  +2.29%  "
    "D14" = "9.72"
    "E14" = "  +3.71%  "
    "D15" = "4.171.56"
    "E15" = "  +4.14%  "
    "D16" = "3.600.37"
    "E16" = "  +3.97%  "
    "D17" = "19.31"
    "E17" = "  +4.17%  "
    "D18" = "69.895.03"
    "E18" = "  +3.95%  "
    "D19" = "12.40"
    "E19" = "  +3.32%  "
    "D20" = "0.120"
    "E20" = "  +0.23%  "
    "D21" = "1.04"
    "E21" = "  +3.12%  "
    "D22" = "488.27"
    "E22" = "  +1.51%  "
    "D23" = "17.38"
    "E23" = "  +15.12%  "
    "D24" = "5.36"
    "E24" = "  +8.10%  "
    "D25" = "4.42"
    "E25" = "  +5.80%  "
    "D26" = "90.16"
    "E26" = "  +0.45%  "
    "D27" = "3.08"
    "E27" = "  +4.15%  "
    "D28" = "11.00"
    "E28" = "  +0.82%  "
    "D29" = "9.35"
    "E29" = "  +4.78%  "
    "D30" = "32.15"
    "D31" = "7.44"
    "E31" = "  +6.47%  "
    "D32" = "622.23"
    "E32" = "  +3.38%  "
    "D33" = "12.18"
    "E33" = "  +4.90%  "
    "E34" = "  +6.25%  "
    "E35" = "  +3.27%  "
    "E36" = "  +3.82%  "
    "D37" = "0.999"
    "E37" = "  +0.02%  "
    "D38" = "0.401"
    "E38" = "  +3.31%  "
    "D39" = "37.74"
    "E39" = "  +3.31%  "
    "E40" = "  -1.20%  "
    "D41" = "3.60"
    "E41" = "  -1.77%  "
    "D42" = "3.295.70"
    "E42" = "  +4.56%  "
    "D43" = "3.07"
    "E43" = "  +4.99%  "
    "E44" = "  +3.88%  "
    "E46" = "  +2.45%  "
    "D47" = "0.136"
    "E47" = "  +1.06%  "
    "D48" = "9.05"
    "E48" = "  +3.51%  "
    "E49" = "  +5.72%  "
    "D50" = "2.67"
    "E50" = "  -4.95%  "
    "E51" = "  +0.20%  "
}

# Cells whose new text would otherwise be auto-parsed as a number by
# Excel (e.g. "1.04", "0.999"); force text format first so the exact
# original string is preserved, matching the source inline-string cells.
$textForceCells = @(
    "D5"
    "D6"
    "D7"
    "D11"
    "D12"
    "D13"
    "D14"
    "D17"
    "D19"
    "D20"
    "D21"
    "D22"
    "D23"
    "D24"
    "D25"
    "D26"
    "D27"
    "D28"
    "D29"
    "D30"
    "D31"
    "D32"
    "D33"
    "D37"
    "D38"
    "D39"
    "D41"
    "D43"
    "D47"
    "D48"
    "D50"
)

foreach ($addr in $cellUpdates.Keys) {
    $cell = $ws.Range($addr)
    if ($textForceCells -contains $addr) {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $cellUpdates[$addr]
}
